# Adding code for new call functionality
$wb = $excel.ActiveWorkbook

# --- Update selection on the existing "Cases" sheet (was the active tab) ---
$casesSheet = $wb.Worksheets("Cases")
$casesSheet.Rows("1:1").Select()

# --- Add the new "Calls" sheet after "Cases" (becomes the active sheet/tab) ---
$callsSheet = $wb.Worksheets.Add([System.Type]::Missing, $casesSheet)
$callsSheet.Name = "Calls"

# Header row (left to right)
$callsSheet.Range("A1").Value = "contact"
$callsSheet.Range("B1").Value = "deal"
$callsSheet.Range("C1").Value = "task"
$callsSheet.Range("D1").Value = "case"
$callsSheet.Range("E1").Value = "notes"
$callsSheet.Range("A1:E1").Interior.Color = 65535

# Data, column by column (top to bottom)
$callsSheet.Range("A2").Value = "zzzx"
$callsSheet.Range("A3").Value = "aaaa"

$callsSheet.Range("B2").Value = "aaaa"
$callsSheet.Range("B3").Value = "bbbb"

$callsSheet.Range("C2").Value = "cccc"
$callsSheet.Range("C3").Value = "dddd"

$callsSheet.Range("D2").Value = "eeee"
$callsSheet.Range("D3").Value = "ffff"

$callsSheet.Range("E2").Value = "gggg"
$callsSheet.Range("E3").Value = "hhhh"

# Selection on the new sheet: full column B selected, active cell B1
$callsSheet.Columns("B:B").Select()
